# "add mysql connect info"
#
# The "Property" sheet lists config-struct fields (Id / Type / Public / Private /
# Save / View / Index / SaveInterval / RelationValue / Desc). Row 5 used to be the
# lone "Pwd" field; this change turns it into a small "SqlXXX" block describing the
# SQL connection: SqlIP (renamed from Pwd, reusing row 5), SqlPort, SqlUser and
# SqlPwd (three brand-new rows appended after it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: "Pwd" -> "SqlIP" -------------------------------------------------
# Only A5 referenced this shared string, so renaming it in place keeps every
# other row's data untouched.
$ws.Range("A5").Value = "SqlIP"

# A5 previously had no explicit style; the updated row gives it the same
# style as the other "key" cells in column A for rows 3/4 (style index 2,
# carried by B5/I5/J5 already). Grab it from B5 via a format-only paste so
# A5's value/type are left alone.
$ws.Range("B5").Copy()
$ws.Range("A5").PasteSpecial(-4122)   # xlPasteFormats

# --- New rows 6-8: SqlPort / SqlUser / SqlPwd --------------------------------
# Same shape as the existing rows: Type, four bool flags (Public/Private/
# Save/View), Index, SaveInterval (both 0), RelationValue=Friend, Desc blank.

$ws.Range("A6").Value = "SqlPort"
$ws.Range("B6").Value = "int"
$ws.Range("C6").Value = $true
$ws.Range("D6").Value = $false
$ws.Range("E6").Value = $false
$ws.Range("F6").Value = $true
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = "Friend"

$ws.Range("A7").Value = "SqlUser"
$ws.Range("B7").Value = "string"
$ws.Range("C7").Value = $true
$ws.Range("D7").Value = $false
$ws.Range("E7").Value = $false
$ws.Range("F7").Value = $true
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = "Friend"

$ws.Range("A8").Value = "SqlPwd"
$ws.Range("B8").Value = "string"
$ws.Range("C8").Value = $true
$ws.Range("D8").Value = $false
$ws.Range("E8").Value = $false
$ws.Range("F8").Value = $true
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = "Friend"

# Carry the row-5 formatting (columns B:J) down onto the three new rows -
# string/int cells use style 2, the boolean/index cells use style 3, same as
# every other data row in this table. Column A is left at the default style
# for rows 7/8 (matching the un-styled A cells the author typed there), but
# row 6's A cell picks up style 2 like A5 above.
$ws.Range("B5:J5").Copy()
$ws.Range("B6:J6").PasteSpecial(-4122)
$ws.Range("B7:J7").PasteSpecial(-4122)
$ws.Range("B8:J8").PasteSpecial(-4122)

$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)

# --- Selection ---------------------------------------------------------------
$ws.Range("C15").Select() | Out-Null
